$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to be treated as literal text so numeric-looking
# strings (e.g. "381.80") are not auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = '51.606.30'
$ws.Cells.Item(2, 5).Value = '  +1.18%  '
$ws.Cells.Item(3, 4).Value = '2.999.30'
$ws.Cells.Item(3, 5).Value = '  +2.18%  '
$ws.Cells.Item(4, 4).Value = '0.999'
$ws.Cells.Item(4, 5).Value = '  -0.13%  '
$ws.Cells.Item(5, 4).Value = '381.80'
$ws.Cells.Item(5, 5).Value = '  +2.10%  '
$ws.Cells.Item(6, 4).Value = '103.48'
$ws.Cells.Item(6, 5).Value = '  +2.59%  '
$ws.Cells.Item(7, 5).Value = '  +2.29%  '
$ws.Cells.Item(8, 5).Value = '  -0.03%  '
$ws.Cells.Item(9, 5).Value = '  +2.41%  '
$ws.Cells.Item(10, 4).Value = '36.88'
$ws.Cells.Item(10, 5).Value = '  +1.96%  '
$ws.Cells.Item(11, 5).Value = '  -0.60%  '
$ws.Cells.Item(12, 5).Value = '  +1.31%  '
$ws.Cells.Item(13, 4).Value = '3.470.48'
$ws.Cells.Item(13, 5).Value = '  +2.09%  '
$ws.Cells.Item(14, 5).Value = '  +4.19%  '
$ws.Cells.Item(15, 4).Value = '18.46'
$ws.Cells.Item(15, 5).Value = '  +2.78%  '
$ws.Cells.Item(16, 4).Value = '3.012.27'
$ws.Cells.Item(16, 5).Value = '  +2.44%  '
$ws.Cells.Item(17, 4).Value = '11.21'
$ws.Cells.Item(17, 5).Value = '  +4.95%  '
$ws.Cells.Item(18, 5).Value = '  +2.59%  '
$ws.Cells.Item(19, 4).Value = '51.611.22'
$ws.Cells.Item(19, 5).Value = '  +1.16%  '
$ws.Cells.Item(20, 4).Value = '3.11'
$ws.Cells.Item(20, 5).Value = '  +0.12%  '
$ws.Cells.Item(21, 4).Value = '12.62'
$ws.Cells.Item(21, 5).Value = '  +1.89%  '
$ws.Cells.Item(22, 4).Value = '0.0₃0966'
$ws.Cells.Item(22, 5).Value = '  +1.07%  '
$ws.Cells.Item(23, 4).Value = '70.52'
$ws.Cells.Item(23, 5).Value = '  +2.81%  '
$ws.Cells.Item(24, 4).Value = '268.38'
$ws.Cells.Item(24, 5).Value = '  +1.26%  '
$ws.Cells.Item(25, 5).Value = '  +3.93%  '
$ws.Cells.Item(26, 4).Value = '7.91'
$ws.Cells.Item(26, 5).Value = '  -1.89%  '
$ws.Cells.Item(27, 4).Value = '7.37'
$ws.Cells.Item(27, 5).Value = '  -3.03%  '
$ws.Cells.Item(28, 5).Value = '  -0.03%  '
$ws.Cells.Item(29, 4).Value = '26.12'
$ws.Cells.Item(29, 5).Value = '  +2.15%  '
$ws.Cells.Item(30, 4).Value = '0.167'
$ws.Cells.Item(30, 5).Value = '  +2.14%  '
$ws.Cells.Item(31, 5).Value = '  -0.17%  '
$ws.Cells.Item(32, 4).Value = '10.37'
$ws.Cells.Item(32, 5).Value = '  +4.23%  '
$ws.Cells.Item(33, 4).Value = '34.68'
$ws.Cells.Item(33, 5).Value = '  +4.30%  '
$ws.Cells.Item(34, 4).Value = '51.58'
$ws.Cells.Item(34, 5).Value = '  +1.74%  '
$ws.Cells.Item(35, 4).Value = '2.07'
$ws.Cells.Item(35, 5).Value = '  +1.34%  '
$ws.Cells.Item(36, 5).Value = '  +0.27%  '
$ws.Cells.Item(37, 5).Value = '  -0.13%  '
$ws.Cells.Item(38, 5).Value = '  +3.70%  '
$ws.Cells.Item(39, 4).Value = '16.82'
$ws.Cells.Item(39, 5).Value = '  +3.59%  '
$ws.Cells.Item(40, 5).Value = '  +2.12%  '
$ws.Cells.Item(41, 5).Value = '  +3.55%  '
$ws.Cells.Item(42, 4).Value = '1.85'
$ws.Cells.Item(42, 5).Value = '  +3.67%  '
$ws.Cells.Item(43, 4).Value = '125.32'
$ws.Cells.Item(43, 5).Value = '  +4.45%  '
$ws.Cells.Item(44, 4).Value = '3.66'
$ws.Cells.Item(44, 5).Value = '  +10.28%  '
$ws.Cells.Item(45, 4).Value = '21.65'
$ws.Cells.Item(45, 5).Value = '  +1.97%  '
$ws.Cells.Item(46, 2).Value = 'ApeXProtocol'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Cells.Item(46, 4).Value = '2.40'
$ws.Cells.Item(46, 5).Value = '  +4.10%  '
$ws.Cells.Item(47, 2).Value = 'WEMIXToken'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Cells.Item(47, 4).Value = '2.03'
$ws.Cells.Item(47, 5).Value = '  +0.02%  '
$ws.Cells.Item(48, 5).Value = '  +0.38%  '
$ws.Cells.Item(49, 4).Value = '2.049.76'
$ws.Cells.Item(49, 5).Value = '  +2.76%  '
$ws.Cells.Item(50, 4).Value = '0.0337'
$ws.Cells.Item(50, 5).Value = '  +4.20%  '
$ws.Cells.Item(51, 4).Value = '0.542'
$ws.Cells.Item(51, 5).Value = '  +17.86%  '

# Restore the original (default/general) formatting on the Price column
# now that the literal text values are safely stored.
$ws.Range("D2:D51").ClearFormats()
